$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.091.34'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -0.70%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.357.80'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -0.93%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.08%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '570.37'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.75%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '135.62'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.35%  '

# Row 7
$ws.Range('E7').Value = '  -0.05%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.356.00'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.95%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.468'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -1.90%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.47'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.13%  '

# Row 11
$ws.Range('E11').Value = '  -3.02%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.385'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -2.37%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.928.48'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -1.04%  '

# Row 14
$ws.Range('E14').Value = '  +1.29%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '25.93'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +2.49%  '

# Row 16
$ws.Range('E16').Value = '  -4.09%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.355.38'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.01%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '61.200.90'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.71%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.99'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.49%  '

# Row 20
$ws.Range('E20').Value = '  -1.27%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.24'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -1.57%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '377.23'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -2.79%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.552'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -3.38%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.486.88'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -1.11%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.07%  '

# Row 26
$ws.Range('B26').Value = 'PEPE'
$ws.Range('C26').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000124'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -2.06%  '

# Row 27
$ws.Range('B27').Value = 'Litecoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '70.91'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.24%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.78'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +11.70%  '

# Row 29
$ws.Range('E29').Value = '  +0.37%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.43'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -3.17%  '

# Row 31
$ws.Range('E31').Value = '  +4.55%  '

# Row 32
$ws.Range('E32').Value = '  -2.12%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.13'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.16%  '

# Row 34
$ws.Range('E34').Value = '  -0.03%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '23.54'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +0.27%  '

# Row 36
$ws.Range('E36').Value = '  -5.53%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.76'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -2.88%  '

# Row 38
$ws.Range('E38').Value = '  -1.43%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '164.68'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +1.72%  '

# Row 40
$ws.Range('E40').Value = '  -5.03%  '

# Row 41
$ws.Range('E41').Value = '  -0.09%  '

# Row 42
$ws.Range('E42').Value = '  -1.27%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.767'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.64%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '41.45'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.31%  '

# Row 45
$ws.Range('E45').Value = '  -1.22%  '

# Row 46
$ws.Range('E46').Value = '  -1.50%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '23.72'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -4.24%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '23.12'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.84%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.79'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -2.49%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.347.13'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.99%  '

# Row 51
$ws.Range('E51').Value = '  -1.94%  '
